$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10, shifting existing rows 10-15 down to 11-16
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record
$ws.Cells.Item(10, 1).Value = "无锡红豆"
$ws.Cells.Item(10, 2).Value = 10000000
$ws.Cells.Item(10, 3).NumberFormat = "@"
$ws.Cells.Item(10, 3).Value = "2024-01-02"
$ws.Cells.Item(10, 4).Value = "瑞e保"
$ws.Cells.Item(10, 5).Value = "无锡红豆居家服饰有限公司"
$ws.Cells.Item(10, 6).Value = "无锡红豆"
